# Trade #5 (row index, "Trade #" value 33) closed at 2026-02-18 00:09:35 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook to reflect Trade #33
# (MarketMaking strategy) transitioning from OPEN to CLOSED, plus the
# knock-on aggregate stats on the Summary and Strategy Status sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - refresh headline stats now that another trade closed
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.91   # Current Capital
$summary.Range("B4").Value = 1.01      # Total P&L $
$summary.Range("B6").Value = 33        # Total Trades
$summary.Range("B7").Value = 18        # Winning Trades
$summary.Range("B9").Value = 54.55     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.91
$status.Range("D6").Value = 4
$status.Range("E6").Value = 0.1
$status.Range("F6").Value = -0.09
$status.Range("G6").Value = 50

# ---------------------------------------------------------------------
# All Trades sheet - Trade #33 row (sheet row 34)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G34").Value = 0.382353
$allTrades.Range("H34").Value = "CLOSED"
$allTrades.Range("I34").Value = 9.2437
$allTrades.Range("J34").Value = 0.03
$allTrades.Range("K34").Value = 99.91
$allTrades.Range("L34").Value = "early_exit"
$allTrades.Range("M34").Value = 0.19

# ---------------------------------------------------------------------
# MarketMaking sheet - same trade, strategy-specific log (row 5)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G5").Value = 0.382353
$marketMaking.Range("H5").Value = "CLOSED"
$marketMaking.Range("I5").Value = 9.2437
$marketMaking.Range("J5").Value = 0.03
$marketMaking.Range("K5").Value = 99.91
$marketMaking.Range("P5").Value = "early_exit"
$marketMaking.Range("Q5").Value = 0.19
